$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '62.891.12'
Set-TextValue $ws 'E2' '  -2.30%  '
Set-TextValue $ws 'D3' '3.390.06'
Set-TextValue $ws 'E3' '  -3.31%  '
Set-TextValue $ws 'D4' '0.999'
Set-TextValue $ws 'E4' '  -0.07%  '
Set-TextValue $ws 'D5' '575.13'
Set-TextValue $ws 'E5' '  -2.72%  '
Set-TextValue $ws 'D6' '126.25'
Set-TextValue $ws 'E7' '  -0.01%  '
Set-TextValue $ws 'D8' '3.391.90'
Set-TextValue $ws 'E8' '  -3.23%  '
Set-TextValue $ws 'E9' '  -2.17%  '
Set-TextValue $ws 'D10' '7.29'
Set-TextValue $ws 'E10' '  -4.47%  '
Set-TextValue $ws 'E11' '  -4.59%  '
Set-TextValue $ws 'E12' '  -3.40%  '
Set-TextValue $ws 'D13' '3.967.71'
Set-TextValue $ws 'E13' '  -3.34%  '
Set-TextValue $ws 'E14' '  -1.09%  '
Set-TextValue $ws 'D15' '3.391.41'
Set-TextValue $ws 'E15' '  -3.01%  '
Set-TextValue $ws 'E16' '  -5.49%  '
Set-TextValue $ws 'D17' '62.907.90'
Set-TextValue $ws 'E17' '  -2.26%  '
Set-TextValue $ws 'D18' '24.75'
Set-TextValue $ws 'E18' '  -4.14%  '
Set-TextValue $ws 'D19' '9.26'
Set-TextValue $ws 'E19' '  -8.02%  '
Set-TextValue $ws 'E20' '  -2.29%  '
Set-TextValue $ws 'D21' '13.15'
Set-TextValue $ws 'E21' '  -3.53%  '
Set-TextValue $ws 'D22' '370.64'
Set-TextValue $ws 'E22' '  -5.64%  '
Set-TextValue $ws 'D24' '3.524.00'
Set-TextValue $ws 'E24' '  -3.36%  '
Set-TextValue $ws 'D25' '1.00'
Set-TextValue $ws 'E25' '  +0.07%  '
Set-TextValue $ws 'D26' '71.74'
Set-TextValue $ws 'E26' '  -3.65%  '
Set-TextValue $ws 'E27' '  -9.52%  '
Set-TextValue $ws 'D28' '0.998'
Set-TextValue $ws 'E28' '  -2.05%  '
Set-TextValue $ws 'D29' '7.04'
Set-TextValue $ws 'E29' '  -5.45%  '
Set-TextValue $ws 'B30' 'InternetComputer(DFINITY)'
Set-TextValue $ws 'C30' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws 'D30' '7.86'
Set-TextValue $ws 'E30' '  -4.24%  '
Set-TextValue $ws 'B31' 'Fetch.AI'
Set-TextValue $ws 'C31' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws 'D31' '1.42'
Set-TextValue $ws 'E31' '  -4.02%  '
Set-TextValue $ws 'B32' 'PancakeSwap'
Set-TextValue $ws 'C32' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws 'D32' '2.12'
Set-TextValue $ws 'E32' '  -6.99%  '
Set-TextValue $ws 'E33' '  -0.01%  '
Set-TextValue $ws 'E34' '  -5.39%  '
Set-TextValue $ws 'D35' '3.417.74'
Set-TextValue $ws 'E35' '  -3.31%  '
Set-TextValue $ws 'D36' '22.72'
Set-TextValue $ws 'E36' '  -2.95%  '
Set-TextValue $ws 'D37' '5.40'
Set-TextValue $ws 'E37' '  +0.87%  '
Set-TextValue $ws 'D38' '166.37'
Set-TextValue $ws 'E38' '  +0.05%  '
Set-TextValue $ws 'D39' '6.67'
Set-TextValue $ws 'E39' '  -4.27%  '
Set-TextValue $ws 'D40' '1.50'
Set-TextValue $ws 'E40' '  -4.58%  '
Set-TextValue $ws 'D41' '0.0756'
Set-TextValue $ws 'E41' '  -4.12%  '
Set-TextValue $ws 'D42' '0.999'
Set-TextValue $ws 'E42' '  -0.12%  '
Set-TextValue $ws 'D43' '41.90'
Set-TextValue $ws 'E43' '  -0.34%  '
Set-TextValue $ws 'D44' '0.764'
Set-TextValue $ws 'E44' '  -5.80%  '
Set-TextValue $ws 'E45' '  -4.72%  '
Set-TextValue $ws 'E46' '  -6.54%  '
Set-TextValue $ws 'D47' '1.11'
Set-TextValue $ws 'E47' '  -6.42%  '
Set-TextValue $ws 'D48' '22.50'
Set-TextValue $ws 'E48' '  -10.22%  '
Set-TextValue $ws 'E49' '  -2.99%  '
Set-TextValue $ws 'D50' '2.240.76'
Set-TextValue $ws 'E50' '  -6.27%  '
Set-TextValue $ws 'D51' '0.843'
Set-TextValue $ws 'E51' '  -8.30%  '
